$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item description")
$ws.Range("C1").EntireColumn.Delete()
